$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 840.9375
$ws.Range("I2").Value = 771.25
$ws.Range("K2").Value = 771.25
$ws.Range("M2").Value = -658.25

$ws.Range("H33").Value = 192.72728
$ws.Range("I33").Value = 127
$ws.Range("K33").Value = 127
$ws.Range("M33").Value = 102

$ws.Range("H40").Value = 3355.1428
$ws.Range("I40").Value = 3581
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 3581
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -3406
$ws.Range("N40").Value = -2350

$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 36000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -40920

$ws.Range("H138").Value = 1532.8
$ws.Range("I138").Value = 1118.6428
$ws.Range("K138").Value = 3355.9284
$ws.Range("M138").Value = 1784.0716

$ws.Range("H139").Value = 70766.71000000001
$ws.Range("J139").Value = 70766.71000000001
$ws.Range("L139").Value = 70766.71000000001
$ws.Range("N139").Value = -81046.71000000001

$ws.Range("H141").Value = 1954.7097
$ws.Range("J141").Value = 5360
$ws.Range("L141").Value = 16080
$ws.Range("N141").Value = -26440


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11278.595
$ws.Range("I32").Value = 5408.905
$ws.Range("J32").Value = 18982.562
$ws.Range("K32").Value = 5408.905
$ws.Range("L32").Value = 18982.562
$ws.Range("M32").Value = -5121.905
$ws.Range("N32").Value = -19556.562

$ws.Range("H51").Value = 5047
$ws.Range("J51").Value = 5047
$ws.Range("L51").Value = 5047
$ws.Range("N51").Value = -6559

$ws.Range("H61").Value = 127976.125
$ws.Range("I61").Value = 3401.4285
$ws.Range("K61").Value = 3401.4285
$ws.Range("M61").Value = -3189.4285

$ws.Range("H74").Value = 35685
$ws.Range("I74").Value = 54164.79
$ws.Range("J74").Value = 3765.3635
$ws.Range("K74").Value = 54164.79
$ws.Range("L74").Value = 3765.3635
$ws.Range("M74").Value = -53290.79
$ws.Range("N74").Value = -5513.363499999999

$ws.Range("H77").Value = 35685
$ws.Range("I77").Value = 54164.79
$ws.Range("J77").Value = 3765.3635
$ws.Range("K77").Value = 270823.95
$ws.Range("L77").Value = 18826.8175
$ws.Range("M77").Value = -266455.95
$ws.Range("N77").Value = -27562.8175

$ws.Range("H122").Value = 2276.3845
$ws.Range("I122").Value = 1299.4286
$ws.Range("K122").Value = 3898.2858
$ws.Range("M122").Value = -1448.2858

$ws.Range("H132").Value = 2296.9614
$ws.Range("I132").Value = 2229.261
$ws.Range("K132").Value = 6687.782999999999
$ws.Range("M132").Value = -4157.782999999999

$ws.Range("H136").Value = 127976.125
$ws.Range("I136").Value = 3401.4285
$ws.Range("K136").Value = 10204.2855
$ws.Range("M136").Value = -7654.2855


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 126699.16
$ws.Range("I20").Value = 163669
$ws.Range("K20").Value = 163669
$ws.Range("M20").Value = -163422


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3055.8386
$ws.Range("I31").Value = 2198.1667
$ws.Range("J31").Value = 4243.385
$ws.Range("K31").Value = 2198.1667
$ws.Range("L31").Value = 4243.385
$ws.Range("M31").Value = -1903.1667
$ws.Range("N31").Value = -4833.385

$ws.Range("H33").Value = 5179.857
$ws.Range("I33").Value = 1565.25
$ws.Range("J33").Value = 9999.333000000001
$ws.Range("K33").Value = 1565.25
$ws.Range("L33").Value = 9999.333000000001
$ws.Range("M33").Value = -1186.25
$ws.Range("N33").Value = -10757.333

$ws.Range("H34").Value = 3055.8386
$ws.Range("I34").Value = 2198.1667
$ws.Range("J34").Value = 4243.385
$ws.Range("K34").Value = 2198.1667
$ws.Range("L34").Value = 4243.385
$ws.Range("M34").Value = -1996.1667
$ws.Range("N34").Value = -4647.385

$ws.Range("H122").Value = 3119.4167
$ws.Range("J122").Value = 2966.375
$ws.Range("L122").Value = 8899.125
$ws.Range("N122").Value = -13799.125

$ws.Range("H132").Value = 4177333.8
$ws.Range("I132").Value = 5052856.5
$ws.Range("K132").Value = 15158569.5
$ws.Range("M132").Value = -15156039.5

$ws.Range("H134").Value = 3451405.2
$ws.Range("I134").Value = 3761909
$ws.Range("K134").Value = 11285727
$ws.Range("M134").Value = -11283192


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 89.42856999999999
$ws.Range("J2").Value = 37.5
$ws.Range("L2").Value = 225
$ws.Range("N2").Value = -451

$ws.Range("H5").Value = 700.7308
$ws.Range("I5").Value = 589.8182
$ws.Range("J5").Value = 782.06665
$ws.Range("K5").Value = 1769.4546
$ws.Range("L5").Value = 2346.19995
$ws.Range("M5").Value = -1657.4546
$ws.Range("N5").Value = -2570.19995

$ws.Range("H6").Value = 200339.4
$ws.Range("I6").Value = 200339.4
$ws.Range("K6").Value = 601018.2
$ws.Range("M6").Value = -600905.2

$ws.Range("H7").Value = 4381.4
$ws.Range("I7").Value = 247.5
$ws.Range("J7").Value = 7137.3335
$ws.Range("K7").Value = 742.5
$ws.Range("L7").Value = 21412.0005
$ws.Range("M7").Value = -630.5
$ws.Range("N7").Value = -21636.0005

$ws.Range("H132").Value = 5969.45
$ws.Range("I132").Value = 2081.6
$ws.Range("K132").Value = 18734.4
$ws.Range("M132").Value = -16204.4

$ws.Range("H135").Value = 700.7308
$ws.Range("I135").Value = 589.8182
$ws.Range("J135").Value = 782.06665
$ws.Range("K135").Value = 5308.3638
$ws.Range("L135").Value = 7038.59985
$ws.Range("M135").Value = -2773.3638
$ws.Range("N135").Value = -12108.59985

$ws.Range("H137").Value = 3626.0908
$ws.Range("I137").Value = 2554.111
$ws.Range("K137").Value = 7662.333
$ws.Range("M137").Value = -2562.333


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 18602.2
$ws.Range("J26").Value = 18602.2
$ws.Range("L26").Value = 18602.2
$ws.Range("N26").Value = -19162.2

$ws.Range("H50").Value = 18602.2
$ws.Range("J50").Value = 18602.2
$ws.Range("L50").Value = 18602.2
$ws.Range("N50").Value = -19598.2

$ws.Range("H52").Value = 19992.5
$ws.Range("J52").Value = 19990
$ws.Range("L52").Value = 19990
$ws.Range("N52").Value = -20508

$ws.Range("H57").Value = 19725
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31640

$ws.Range("H58").Value = 28994
$ws.Range("J58").Value = 28742.5
$ws.Range("L58").Value = 28742.5
$ws.Range("N58").Value = -29296.5

$ws.Range("H93").Value = 18215
$ws.Range("J93").Value = 18215
$ws.Range("L93").Value = 18215
$ws.Range("N93").Value = -21959

$ws.Range("H102").Value = 1140
$ws.Range("I102").Value = 1092.6
$ws.Range("K102").Value = 1092.6
$ws.Range("M102").Value = 529.4000000000001

$ws.Range("H117").Value = 53334
$ws.Range("J117").Value = 53334
$ws.Range("L117").Value = 53334
$ws.Range("N117").Value = -60218

$ws.Range("H122").Value = 9354705
$ws.Range("I122").Value = 10204792
$ws.Range("K122").Value = 30614376
$ws.Range("M122").Value = -30611926

$ws.Range("H132").Value = 4431.35
$ws.Range("I132").Value = 3441.8
$ws.Range("J132").Value = 7400
$ws.Range("K132").Value = 10325.4
$ws.Range("L132").Value = 22200
$ws.Range("M132").Value = -7795.400000000001
$ws.Range("N132").Value = -27260


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5979.6
$ws.Range("I7").Value = 4999
$ws.Range("K7").Value = 4999
$ws.Range("M7").Value = -4887

$ws.Range("H16").Value = 1330.3914
$ws.Range("I16").Value = 1356.3684
$ws.Range("J16").Value = 1207
$ws.Range("K16").Value = 1356.3684
$ws.Range("L16").Value = 1207
$ws.Range("M16").Value = -1186.3684
$ws.Range("N16").Value = -1547

$ws.Range("H40").Value = 14145606
$ws.Range("I40").Value = 4223.4
$ws.Range("K40").Value = 4223.4
$ws.Range("M40").Value = -4087.4

$ws.Range("H45").Value = 13697.833
$ws.Range("I45").Value = 14010.25
$ws.Range("J45").Value = 13073
$ws.Range("K45").Value = 14010.25
$ws.Range("L45").Value = 13073
$ws.Range("M45").Value = -13603.25
$ws.Range("N45").Value = -13887

$ws.Range("H46").Value = 1677.25
$ws.Range("I46").Value = 1666.091
$ws.Range("J46").Value = 1800
$ws.Range("K46").Value = 1666.091
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -1478.091
$ws.Range("N46").Value = -2176

$ws.Range("H55").Value = 5263842
$ws.Range("I55").Value = 635
$ws.Range("J55").Value = 7408111.5
$ws.Range("K55").Value = 635
$ws.Range("L55").Value = 7408111.5
$ws.Range("M55").Value = -462
$ws.Range("N55").Value = -7408457.5

$ws.Range("H126").Value = 5979.6
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527

$ws.Range("H132").Value = 1901
$ws.Range("I132").Value = 1901
$ws.Range("K132").Value = 5703
$ws.Range("M132").Value = -3173


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1056.8
$ws.Range("I81").Value = 873.75
$ws.Range("K81").Value = 1747.5
$ws.Range("M81").Value = -686.5

$ws.Range("H84").Value = 1056.8
$ws.Range("I84").Value = 873.75
$ws.Range("K84").Value = 8737.5
$ws.Range("M84").Value = -3433.5

$ws.Range("H132").Value = 2377.1875
$ws.Range("I132").Value = 2214.0454
$ws.Range("K132").Value = 6642.1362
$ws.Range("M132").Value = -4112.1362

